$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 5 (Choose Product -> Form Filling "Untuk ..." subsections):
# prefix validation-scenario labels with "Enter " ---
$ws.Range("B63").Value  = "Enter Nomor Handphone VALID"
$ws.Range("B64").Value  = "Enter Nomor Handphone INVALID"
$ws.Range("B67").Value  = "Enter Nomor Peserta VALID"
$ws.Range("B68").Value  = "Enter Nomor Peserta INVALID"
$ws.Range("B71").Value  = "Enter ID Pelanggan VALID DAN Nomor Handphone VALID"
$ws.Range("B72").Value  = "Enter ID Pelanggan VALID DAN Nomor Handphone INVALID"
$ws.Range("B73").Value  = "Enter ID Pelanggan INVALID DAN Nomor Handphone VALID"
$ws.Range("B74").Value  = "Enter ID Pelanggan INVALID DAN Nomor Handphone INVALID"
$ws.Range("B77").Value  = "Enter Nomor Pelanggan VALID DAN Wilayah VALID"
$ws.Range("B78").Value  = "Enter Nomor Pelanggan VALID DAN Wilayah INVALID"
$ws.Range("B79").Value  = "Enter Nomor Pelanggan INVALID DAN Wilayah VALID"
$ws.Range("B80").Value  = "Enter Nomor Pelanggan INVALID DAN Wilayah INVALID"
$ws.Range("B83").Value  = "Enter ID Pelanggan VALID"
$ws.Range("B84").Value  = "Enter ID Pelanggan INVALID"
$ws.Range("B87").Value  = "Enter Nomor Kartu Kredit VALID"
$ws.Range("B88").Value  = "Enter Nomor Kartu Kredit INVALID"
$ws.Range("B91").Value  = "Enter Kota VALID DAN Tahun Pajak VALID DAN Nomor Objek Pajak VALID"
$ws.Range("B92").Value  = "Enter Kota VALID DAN Tahun Pajak VALID DAN Nomor Objek Pajak INVALID"
$ws.Range("B93").Value  = "Enter Kota VALID DAN Tahun Pajak INVALID DAN Nomor Objek Pajak VALID"
$ws.Range("B94").Value  = "Enter Kota VALID DAN Tahun Pajak INVALID DAN Nomor Objek Pajak INVALID"
$ws.Range("B95").Value  = "Enter Kota INVALID DAN Tahun Pajak VALID DAN Nomor Objek Pajak VALID"
$ws.Range("B96").Value  = "Enter Kota INVALID DAN Tahun Pajak VALID DAN Nomor Objek Pajak INVALID"
$ws.Range("B97").Value  = "Enter Kota INVALID DAN Tahun Pajak INVALID DAN Nomor Objek Pajak VALID"
$ws.Range("B98").Value  = "Enter Kota INVALID DAN Tahun Pajak INVALID DAN Nomor Objek Pajak INVALID"
$ws.Range("B101").Value = "Enter Nama Properti VALID DAN Nomor Pelanggan VALID"
$ws.Range("B102").Value = "Enter Nama Properti VALID DAN Nomor Pelanggan INVALID"
$ws.Range("B103").Value = "Enter Nama Properti INVALID DAN Nomor Pelanggan VALID"
$ws.Range("B104").Value = "Enter Nama Properti INVALID DAN Nomor Pelanggan INVALID"

# --- Section 7 (Payment Method Functionality): rename header, rework / extend
# the Buy-Pulsa payment scenarios (Kredivo + new Gopay cases) ---
$ws.Range("A106").Value = "Buy Pulsa Payment Method Functionality"
$ws.Range("B107").Value = "Enter Nomor Handphone VALID DAN Kredivo"
$ws.Range("B108").Value = "Enter Nomor Handphone INVALID DAN Kredivo"
$ws.Range("B109").Value = "Email Nomor Handphone VALID DAN Gopay"
$ws.Range("B110").Value = "Email Nomor Handphone INVALID DAN Gopay"

# --- restore the cursor/selection position left behind by the author ---
$ws.Range("B39").Select()
